$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "sequence/run_0684_samples/" containing path from each fastq
# filename stored in column F (rows 2 through 19).
$prefix = "sequence/run_0684_samples/"
for ($r = 2; $r -le 19; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value2
    if ($val -and $val.ToString().StartsWith($prefix)) {
        $cell.Value = $val.ToString().Substring($prefix.Length)
    }
}

# Update the active selection to reflect the column of interest (F2:F19)
$ws.Range("F2:F19").Select()
